# g13.1a.xlsx update: refresh quarter data (01/01/2024 -> 01/04/2024),
# swap Distrito Federal/Santa Catarina and Rio Grande do Sul/Mato Grosso rows,
# update Valor figures and Sergipe's Colocação ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the style of a plain data cell (A2) so we can restore it after
# touching NumberFormat on column C -- keeps the text values in C from
# being auto-converted to Excel date serials while leaving cell styling
# untouched.
$plainStyle = $ws.Cells.Item(2, 1).Style

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $plainStyle
}

# --- Região (column A) swaps ---
$ws.Range("A2").Value = "Santa Catarina"
$ws.Range("A3").Value = "Distrito Federal"
$ws.Range("A6").Value = "Mato Grosso"
$ws.Range("A7").Value = "Rio Grande do Sul"

# --- Trimestre (column C): 01/01/2024 -> 01/04/2024 for every data row ---
for ($r = 2; $r -le 10; $r++) {
    Set-TextValue $r 3 "01/04/2024"
}

# --- Valor (column D) updates ---
$ws.Range("D2").Value = 55.95080666490346
$ws.Range("D3").Value = 55.8125
$ws.Range("D4").Value = 55.1427191195312
$ws.Range("D5").Value = 54.70257234726687
$ws.Range("D6").Value = 54.63576158940398
$ws.Range("D7").Value = 53.8840830449827
$ws.Range("D8").Value = 46.69732441471572
$ws.Range("D9").Value = 43.2565479151301
$ws.Range("D10").Value = 50.4425227718206

# --- Colocação (column E) update ---
$ws.Range("E8").Value = "17º"
